# Daily rollover update: for every data row, the "剩余" (days remaining)
# counter is recomputed against "today" (the day this automation runs),
# based on the row's start date (F, stored as a YYYYMMDD integer) plus its
# total-day allotment (D). When the countdown would reach zero or below,
# the cycle restarts: the remaining-days counter resets to the full
# allotment and the start date rolls to today.
#
# Today for this run is 2026-02-16 (the run date implied by the source
# commit, one day after the prior snapshot's 2026-02-15).

function DaysFromYMD($y, $m, $d) {
    # Howard Hinnant's civil_from_days / days_from_civil algorithm,
    # done with plain integer arithmetic (no locale-sensitive date
    # parsing needed since F is just a YYYYMMDD integer).
    if ($m -le 2) { $y = $y - 1 }
    if ($y -ge 0) {
        $era = [math]::Floor($y / 400)
    } else {
        $era = [math]::Floor(($y - 399) / 400)
    }
    $yoe = $y - $era * 400
    if ($m -gt 2) { $mp = $m - 3 } else { $mp = $m + 9 }
    $doy = [math]::Floor((153 * $mp + 2) / 5) + $d - 1
    $doe = $yoe * 365 + [math]::Floor($yoe / 4) - [math]::Floor($yoe / 100) + $doy
    return $era * 146097 + $doe - 719468
}

function DaysFromYMDInt($ymd) {
    $y = [math]::Floor($ymd / 10000)
    $m = [math]::Floor(($ymd % 10000) / 100)
    $d = $ymd % 100
    return DaysFromYMD $y $m $d
}

function YMDIntFromDays($z) {
    $z = $z + 719468
    if ($z -ge 0) {
        $era = [math]::Floor($z / 146097)
    } else {
        $era = [math]::Floor(($z - 146096) / 146097)
    }
    $doe = $z - $era * 146097
    $yoe = [math]::Floor(($doe - [math]::Floor($doe / 1460) + [math]::Floor($doe / 36524) - [math]::Floor($doe / 146096)) / 365)
    $y = $yoe + $era * 400
    $doy = $doe - (365 * $yoe + [math]::Floor($yoe / 4) - [math]::Floor($yoe / 100))
    $mp = [math]::Floor((5 * $doy + 2) / 153)
    $d = $doy - [math]::Floor((153 * $mp + 2) / 5) + 1
    if ($mp -lt 10) { $m = $mp + 3 } else { $m = $mp - 9 }
    if ($m -le 2) { $y = $y + 1 }
    return ($y * 10000) + ($m * 100) + $d
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$todayYmd = 20260216
$todayDays = DaysFromYMDInt $todayYmd

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $total = $ws.Cells.Item($r, 4).Value2
    $remain = $ws.Cells.Item($r, 5).Value2
    $startYmd = $ws.Cells.Item($r, 6).Value2

    if ($total -eq $null -or $startYmd -eq $null) { continue }

    $y = [math]::Floor($startYmd / 10000)
    $m = [math]::Floor(($startYmd % 10000) / 100)
    $d = $startYmd % 100
    if ($y -lt 1900 -or $y -gt 2200 -or $m -lt 1 -or $m -gt 12 -or $d -lt 1 -or $d -gt 31) {
        # Malformed start date (e.g. stray 9-digit value) - leave the row
        # exactly as-is, matching the source automation's behaviour.
        continue
    }

    $startDays = DaysFromYMDInt $startYmd
    $endDays = $startDays + $total
    $newRemain = $endDays - $todayDays

    if ($newRemain -le 0) {
        $ws.Cells.Item($r, 5).Value = $total
        $ws.Cells.Item($r, 6).Value = $todayYmd
    } else {
        $ws.Cells.Item($r, 5).Value = $newRemain
    }
}
